{"js": "// Apply the goods-receipt line-item edit: the single product line in the\n// second table is changed from the \"Admonter\" board to the \"Bolon, Silence\"\n// package, along with its related quantity/unit/price figures.\n\n// 1) Product name (English) - unique text in the document body.\nconst nameEn = context.document.body.search(\n  \"Admonter, alpin, natural oiled, 15x80x2380\",\n  { matchCase: true }\n);\nawait context.sync();\nnameEn.items[0].insertText(\"Bolon, Silence, 2x754x754\", \"Replace\");\n\n// 2) \"\u041c\u0435\u0441\u0442\" (places) column: \"20|\" -> \"0|\"\nconst mest = context.document.body.search(\"20|\", { matchCase: true });\nawait context.sync();\nmest.items[0].insertText(\"0|\", \"Replace\");\n\n// 3) Package type: \"board\" -> \"package\"\nconst pkg = context.document.body.search(\"board\", { matchCase: true });\nawait context.sync();\npkg.items[0].insertText(\"package\", \"Replace\");\n\n// 4) Quantity value (first run of the \"2.14 \u043a\u0433\" cell, keeps \"\u043a\u0433\" run intact):\n//    \"2.14 \" -> \"2.9 \"\nconst qty = context.document.body.search(\"2.14 \", { matchCase: true });\nawait context.sync();\nqty.items[0].insertText(\"2.9 \", \"Replace\");\n\n// 5) Russian product name - unique text in the document body.\nconst nameRu = context.document.body.search(\n  \"\u0410\u0434\u043c\u043e\u043d\u0442\u0435\u0440, \u0431\u0440\u0430\u0448\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u044c\u043f\u0438\u043d, \u043f\u043e\u0434 \u043d\u0430\u0442\u0443\u0440\u0430\u043b\u044c\u043d\u044b\u043c \u043c\u0430\u0441\u043b\u043e\u043c, 15x80x2380\",\n  { matchCase: true }\n);\nawait context.sync();\nnameRu.items[0].insertText(\"\u0411\u043e\u043b\u043e\u043d, Silence, 2x754x754\", \"Replace\");\n\nawait context.sync();\n\n// The remaining three cells (\"5\" -> \"0\", \"\u0448\u0442\" -> \"\u043c2\", \"0\" -> \"121\") are not\n// unique strings in the document, so address them by their table/row/column\n// position instead of by text search. Replace each cell's content via its\n// range (rather than the blunter `cell.value =` setter) so the existing\n// paragraph/run formatting is preserved exactly.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst dataTable = tables.items[1];\ndataTable.rows.load(\"items\");\nawait context.sync();\n\nconst dataRow = dataTable.rows.items[3];\ndataRow.cells.load(\"items\");\nawait context.sync();\n\nconst qtyCountCell = dataRow.cells.items[9]; // \"5\" -> \"0\"\nqtyCountCell.body.getRange(\"Content\").insertText(\"0\", \"Replace\");\n\nconst qtyUnitCell = dataRow.cells.items[10]; // \"\u0448\u0442\" -> \"\u043c2\"\nqtyUnitCell.body.getRange(\"Content\").insertText(\"\u043c2\", \"Replace\");\n\nconst priceCell = dataRow.cells.items[11]; // \"0\" -> \"121\"\npriceCell.body.getRange(\"Content\").insertText(\"121\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Apply the goods-receipt line-item edit: the single product line in the\n# second table is changed from the \"Admonter\" board to the \"Bolon, Silence\"\n# package, along with its related quantity/unit/price figures.\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($findText, $replaceText) {\n    # These strings only occur once each in the document, so a plain\n    # Find/Replace across the whole document body is unambiguous.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $wdFindContinue = 1\n    $wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n\n# 1) Product name (English)\nReplace-UniqueText \"Admonter, alpin, natural oiled, 15x80x2380\" \"Bolon, Silence, 2x754x754\"\n\n# 2) \"\u041c\u0435\u0441\u0442\" (places) column\nReplace-UniqueText \"20|\" \"0|\"\n\n# 3) Package type\nReplace-UniqueText \"board\" \"package\"\n\n# 4) Quantity value (keeps the separate \"\u043a\u0433\" run intact)\nReplace-UniqueText \"2.14 \" \"2.9 \"\n\n# 5) Product name (Russian)\nReplace-UniqueText \"\u0410\u0434\u043c\u043e\u043d\u0442\u0435\u0440, \u0431\u0440\u0430\u0448\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u044c\u043f\u0438\u043d, \u043f\u043e\u0434 \u043d\u0430\u0442\u0443\u0440\u0430\u043b\u044c\u043d\u044b\u043c \u043c\u0430\u0441\u043b\u043e\u043c, 15x80x2380\" \"\u0411\u043e\u043b\u043e\u043d, Silence, 2x754x754\"\n\n# The remaining three cells (\"5\" -> \"0\", \"\u0448\u0442\" -> \"\u043c2\", \"0\" -> \"121\") are not\n# unique strings in the document, so address them by their table/row/column\n# position instead of by text search. Word.Table.Cell is 1-indexed; this is\n# the data row (table 2, row 4) of the line-items table.\n$t = $d.Tables.Item(2)\n$t.Cell(4, 10).Range.Text = \"0\"    # \u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e value: 5 -> 0\n$t.Cell(4, 11).Range.Text = \"\u043c2\"   # \u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e unit: \u0448\u0442 -> \u043c2\n$t.Cell(4, 12).Range.Text = \"121\"  # \u0426\u0435\u043d\u0430: 0 -> 121\n"}
